# Update the F1 fuse BOM row: change from 500mA / 0805 footprint / C2649565
# to 750mA / C328915 / 1206 footprint, per commit "Change to 750mA fuse 1206 package"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "750mA"
$ws.Range("C15").Value = "Fuse_1206_3216Metric_Pad1.42x1.75mm_HandSolder"
$ws.Range("D15").Value = "C328915"

$ws.Range("C15").Select()
